# Commit: add number of comment
# Adds a new "评论数" (number of comments) column in column I of the
# product-info sheet, with a comment-count value for each product row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header
$ws.Range("I1").Value = "评论数"

# Comment-count values per product row (column I)
$ws.Range("I2").Value  = 1
$ws.Range("I3").Value  = 1
$ws.Range("I4").Value  = 1
$ws.Range("I5").Value  = 1
$ws.Range("I6").Value  = 1
$ws.Range("I7").Value  = 1

$ws.Range("I9").Value  = 1
$ws.Range("I10").Value = 1
$ws.Range("I11").Value = 1
$ws.Range("I12").Value = 1
$ws.Range("I13").Value = 1
$ws.Range("I14").Value = 1

$ws.Range("I16").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("I18").Value = 1
$ws.Range("I19").Value = 1
$ws.Range("I20").Value = 1

$ws.Range("I22").Value = 1
$ws.Range("I23").Value = 1

$ws.Range("I25").Value = 0
$ws.Range("I26").Value = 1
$ws.Range("I27").Value = 1

$ws.Range("I29").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("I31").Value = 0

# Leave the view scrolled/selected where data entry finished
$ws.Range("I31").Select()
